$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 978.7692
$ws.Range("I33").Value = 752.7
$ws.Range("K33").Value = 752.7
$ws.Range("M33").Value = -523.7

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 3416
$ws.Range("I40").Value = 2544.6667
$ws.Range("J40").Value = 4984.4
$ws.Range("K40").Value = 2544.6667
$ws.Range("L40").Value = 4984.4
$ws.Range("M40").Value = -2369.6667
$ws.Range("N40").Value = -5334.4

# Row 51: A Bile Business / Shark Oil
$ws.Range("H51").Value = 15632787
$ws.Range("I51").Value = 41671170
$ws.Range("K51").Value = 41671170
$ws.Range("M51").Value = -41670686

# Row 55: A Real Smooth Move / Lanolin
$ws.Range("H55").Value = 20
$ws.Range("I55").Value = 20
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 20
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = 194

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 13528.667
$ws.Range("J137").Value = 5567.1
$ws.Range("L137").Value = 16701.3
$ws.Range("N137").Value = -21801.3

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3713.158
$ws.Range("I138").Value = 3224.0527
$ws.Range("J138").Value = 4202.263
$ws.Range("K138").Value = 9672.158100000001
$ws.Range("L138").Value = 12606.789
$ws.Range("M138").Value = -4532.158100000001
$ws.Range("N138").Value = -22886.789

$ws = $wb.Worksheets.Item("ARM")
# Row 19: Stadium Envy / Bronze Gauntlets
$ws.Range("H19").Value = 9500
$ws.Range("I19").Value = 9000
$ws.Range("K19").Value = 9000
$ws.Range("M19").Value = -8771

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 20273.05
$ws.Range("I32").Value = 22686.334
$ws.Range("J32").Value = 2690.5715
$ws.Range("K32").Value = 22686.334
$ws.Range("L32").Value = 2690.5715
$ws.Range("M32").Value = -22399.334
$ws.Range("N32").Value = -3264.5715

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 3758.0625
$ws.Range("I63").Value = 1945
$ws.Range("K63").Value = 1945
$ws.Range("M63").Value = -1259

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 3758.0625
$ws.Range("I66").Value = 1945
$ws.Range("K66").Value = 9725
$ws.Range("M66").Value = -6293

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2388.8386
$ws.Range("I122").Value = 2313.5789
$ws.Range("J122").Value = 2508
$ws.Range("K122").Value = 6940.736699999999
$ws.Range("L122").Value = 7524
$ws.Range("M122").Value = -4490.736699999999
$ws.Range("N122").Value = -12424

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 1822.8695
$ws.Range("I16").Value = 1492.8572
$ws.Range("K16").Value = 1492.8572
$ws.Range("M16").Value = -1205.8572

# Row 19: Shielding Sales / Square Ash Shield
$ws.Range("H19").Value = 2645.8
$ws.Range("I19").Value = 1082.25
$ws.Range("K19").Value = 1082.25
$ws.Range("M19").Value = -912.25

# Row 24: What You Need / Square Ash Shield
$ws.Range("H24").Value = 2645.8
$ws.Range("I24").Value = 1082.25
$ws.Range("K24").Value = 1082.25
$ws.Range("M24").Value = -912.25

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 7696407.5
$ws.Range("I31").Value = 14287085
$ws.Range("K31").Value = 14287085
$ws.Range("M31").Value = -14286790

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 7696407.5
$ws.Range("I34").Value = 14287085
$ws.Range("K34").Value = 14287085
$ws.Range("M34").Value = -14286883

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 13028.518
$ws.Range("I58").Value = 1914.6
$ws.Range("J58").Value = 37726.11
$ws.Range("K58").Value = 1914.6
$ws.Range("L58").Value = 37726.11
$ws.Range("M58").Value = -1711.6
$ws.Range("N58").Value = -38132.11

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 11348.777
$ws.Range("I62").Value = 13387.9
$ws.Range("J62").Value = 8799.875
$ws.Range("K62").Value = 13387.9
$ws.Range("L62").Value = 8799.875
$ws.Range("M62").Value = -12763.9
$ws.Range("N62").Value = -10047.875

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 11348.777
$ws.Range("I65").Value = 13387.9
$ws.Range("J65").Value = 8799.875
$ws.Range("K65").Value = 66939.5
$ws.Range("L65").Value = 43999.375
$ws.Range("M65").Value = -63819.5
$ws.Range("N65").Value = -50239.375

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 7790.706
$ws.Range("I99").Value = 5118.385
$ws.Range("K99").Value = 5118.385
$ws.Range("M99").Value = -3620.385

# Row 103: Spare a Rod and Spoil the Fishers / Gazelle Horn Fishing Rod
$ws.Range("H103").Value = 20798.4
$ws.Range("I103").Value = 15999.25
$ws.Range("K103").Value = 15999.25
$ws.Range("M103").Value = -14827.25

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 1822.8695
$ws.Range("I113").Value = 1492.8572
$ws.Range("K113").Value = 1492.8572
$ws.Range("M113").Value = 677.1428000000001

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 7790.706
$ws.Range("I126").Value = 5118.385
$ws.Range("K126").Value = 15355.155
$ws.Range("M126").Value = -12885.155

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 33787.453
$ws.Range("I132").Value = 42567.125
$ws.Range("K132").Value = 127701.375
$ws.Range("M132").Value = -125171.375

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 13028.518
$ws.Range("I136").Value = 1914.6
$ws.Range("J136").Value = 37726.11
$ws.Range("K136").Value = 5743.799999999999
$ws.Range("L136").Value = 113178.33
$ws.Range("M136").Value = -3193.799999999999
$ws.Range("N136").Value = -118278.33

$ws = $wb.Worksheets.Item("CUL")
# Row 23: Sweet Smell of Success / Lavender Oil
$ws.Range("H23").Value = 728.8333
$ws.Range("I23").Value = 600
$ws.Range("J23").Value = 740.5454999999999
$ws.Range("K23").Value = 1800
$ws.Range("L23").Value = 2221.6365
$ws.Range("M23").Value = -1565
$ws.Range("N23").Value = -2691.6365

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2408.62
$ws.Range("I132").Value = 2421.0408
$ws.Range("K132").Value = 7263.1224
$ws.Range("M132").Value = -4733.1224

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 926.8
$ws.Range("I22").Value = 724
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = 724
$ws.Range("L22").Value = 1400
$ws.Range("M22").Value = -429
$ws.Range("N22").Value = -1990

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 926.8
$ws.Range("I27").Value = 724
$ws.Range("J27").Value = 1400
$ws.Range("K27").Value = 724
$ws.Range("L27").Value = 1400
$ws.Range("M27").Value = -617
$ws.Range("N27").Value = -1614

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 907.0454999999999
$ws.Range("I55").Value = 372.6
$ws.Range("J55").Value = 1352.4166
$ws.Range("K55").Value = 372.6
$ws.Range("L55").Value = 1352.4166
$ws.Range("M55").Value = -199.6
$ws.Range("N55").Value = -1698.4166

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 1211.8889
$ws.Range("I82").Value = 904.2222
$ws.Range("J82").Value = 1519.5555
$ws.Range("K82").Value = 904.2222
$ws.Range("L82").Value = 1519.5555
$ws.Range("M82").Value = -543.2222
$ws.Range("N82").Value = -2241.5555

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 1211.8889
$ws.Range("I85").Value = 904.2222
$ws.Range("J85").Value = 1519.5555
$ws.Range("K85").Value = 904.2222
$ws.Range("L85").Value = 1519.5555
$ws.Range("M85").Value = 343.7778
$ws.Range("N85").Value = -4015.5555

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 2364.6287
$ws.Range("I132").Value = 1777.5714
$ws.Range("K132").Value = 5332.7142
$ws.Range("M132").Value = -2802.7142

# Row 134: Freezing Fingers / Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value = 47331.168
$ws.Range("J134").Value = 47331.168
$ws.Range("L134").Value = 47331.168
$ws.Range("N134").Value = -57471.168

$ws = $wb.Worksheets.Item("WVR")
# Row 45: Private Concerns / Linen Trousers
$ws.Range("H45").Value = 18541
$ws.Range("I45").Value = 20000
$ws.Range("J45").Value = 17811.5
$ws.Range("K45").Value = 20000
$ws.Range("M45").Value = -19509
$ws.Range("N45").Value = -18793.5
